# Corrige les numeros de groupes : 030502XXXX -> 040311XXXX pour les paroisses.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New group numbers for column C (rows 2..13), replacing 305020100..305021200
$newValues = @(
    403110100,
    403110200,
    403110300,
    403110400,
    403110500,
    403110600,
    403110700,
    403110800,
    403110900,
    403111000,
    403111100,
    403111200
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 3)   # Column C
    $cell.ClearFormats()
    $cell.Value2 = $newValues[$i]
}

# Move the active selection as recorded in the saved view
$ws.Range("E16").Select()
